$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: date, start time, end time (Hours column B left blank)
$ws.Range("A4").Value = 43829

# Match the date formatting already used in A2:A3 (numFmtId 14, "m/d/yyyy")
# by copying the format from A3 instead of re-creating a number format,
# so the same cell style (s="1") is reused.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C4").Value = "1pm"
$ws.Range("D4").Value = "3pm"

# Update selection to match the recorded cursor position after edit
$ws.Range("C6").Select()
